$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 1-10 block: replace several formula cells with plain values ---
$ws.Range("B1").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("B5").Formula = "=24/25"
$ws.Range("B6").Value = 1
$ws.Range("B7").Formula = "=21/25"
$ws.Range("B8").Value = 1
$ws.Range("B9").Value = 1

# --- Row 12 block ---
$ws.Range("B12").Value = 0

# --- Rows 19-40 ---
$ws.Range("B19").Formula = "=7/25"
$ws.Range("B21").Formula = "=6/25"
$ws.Range("B29").Value = 0
$ws.Range("B35").Formula = "=5/25"
$ws.Range("B36").Formula = "=13/25"
$ws.Range("B37").Formula = "=24/25"
$ws.Range("B40").Formula = "=11/25"

# --- Apply Percent style to the average cells ---
$ws.Range("E6").Style = "Percent"
$ws.Range("E17").Style = "Percent"

# --- Selection / active cell ---
[void]$ws.Range("B41").Select()
